# Applies two changes:
#  1. Update the cached "Date Placeholder" field text from 8/7/2018 to 9/21/18
#     across the slide master and every slide layout.
#  2. Fix the spelling of "CrearCommand" -> "ClearCommand" on slide 2.

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            if ($shp.Name -like "Date Placeholder*") {
                if ($shp.TextFrame.TextRange.Text -eq "8/7/2018") {
                    $shp.TextFrame.TextRange.Text = "9/21/18"
                }
            }
        }
    }
}

# Slide master
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every custom (slide) layout
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# Fix "CrearCommand" -> "ClearCommand" spelling on slide 2
$slide2 = $p.Slides.Item(2)
for ($j = 1; $j -le $slide2.Shapes.Count; $j++) {
    $shp = $slide2.Shapes.Item($j)
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.TextRange.Text -eq "CrearCommand") {
            $shp.TextFrame.TextRange.Text = "ClearCommand"
        }
    }
}
